$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing count (column C) values for rows that changed ---
$ws.Range("C960").Value  = 31
$ws.Range("C978").Value  = 39
$ws.Range("C981").Value  = 36
$ws.Range("C985").Value  = 35
$ws.Range("C990").Value  = 16
$ws.Range("C998").Value  = 9
$ws.Range("C1000").Value = 30
$ws.Range("C1010").Value = 38
$ws.Range("C1013").Value = 22
$ws.Range("C1014").Value = 39
$ws.Range("C1019").Value = 39
$ws.Range("C1024").Value = 46
$ws.Range("C1028").Value = 21
$ws.Range("C1029").Value = 38
$ws.Range("C1034").Value = 20
$ws.Range("C1039").Value = 39
$ws.Range("C1043").Value = 29
$ws.Range("C1047").Value = 15
$ws.Range("C1054").Value = 50
$ws.Range("C1057").Value = 21
$ws.Range("C1058").Value = 44
$ws.Range("C1059").Value = 3
$ws.Range("C1062").Value = 23
$ws.Range("C1063").Value = 32
$ws.Range("C1068").Value = 29
$ws.Range("C1069").Value = 39
$ws.Range("C1071").Value = 8
$ws.Range("C1072").Value = 15
$ws.Range("C1073").Value = 34
$ws.Range("C1077").Value = 10
$ws.Range("C1078").Value = 16
$ws.Range("C1079").Value = 29

# --- Rows 1080-1083: the 2020-12-14 (44179) block gains a new 40-49
#     bracket, shifting the existing age-group labels down one bracket
#     and changing their counts ---
$ws.Range("B1080").Value = "40-49"
$ws.Range("C1080").Value = 1

$ws.Range("B1081").Value = "50-59"
$ws.Range("C1081").Value = 4

$ws.Range("B1082").Value = "60-69"
$ws.Range("C1082").Value = 7

$ws.Range("B1083").Value = "70-79"
$ws.Range("C1083").Value = 20

# --- New rows 1084-1088: remainder of the 2020-12-14 block plus the
#     entire new 2020-12-15 (44180) block ---
$ws.Range("A1084").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A1084").Value = 44179
$ws.Range("B1084").Value = "80+"
$ws.Range("C1084").Value = 39

$ws.Range("A1085").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A1085").Value = 44180
$ws.Range("B1085").Value = "50-59"
$ws.Range("C1085").Value = 1

$ws.Range("A1086").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A1086").Value = 44180
$ws.Range("B1086").Value = "60-69"
$ws.Range("C1086").Value = 8

$ws.Range("A1087").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A1087").Value = 44180
$ws.Range("B1087").Value = "70-79"
$ws.Range("C1087").Value = 4

$ws.Range("A1088").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A1088").Value = 44180
$ws.Range("B1088").Value = "80+"
$ws.Range("C1088").Value = 9
